# Lab 3 Rubric update: add a new "PR sent to lab partner" requirement (3 pts)
# and rebalance a few existing point values so the rubric still totals 40.

$wb = $excel.ActiveWorkbook
$wsRubric = $wb.Worksheets.Item("Rubric")
$wsGrade  = $wb.Worksheets.Item("Grade")

# --- Rubric sheet -----------------------------------------------------
# "5 model properties minimum" possible points: 3 -> 2
$wsRubric.Range("B10").Value = 2

# "form for entering msg, story, or post" possible points: 5 -> 4
$wsRubric.Range("B12").Value = 4

# New requirement row: "PR sent to lab partner" worth 3 points
$wsRubric.Range("A19").Value = "PR sent to lab partner"
$wsRubric.Range("B19").Value = 3

# "Site running on Azure" possible points: 5 -> 4
$wsRubric.Range("B20").Value = 4

# --- Grade sheet (mirrors Rubric, with an extra "Actual" column) ------
# "5 model properties minimum": Possible/Actual 3 -> 2
$wsGrade.Range("B11").Value = 2
$wsGrade.Range("C11").Value = 2

# "form for entering msg, story, or post": Possible/Actual 5 -> 4
$wsGrade.Range("B13").Value = 4
$wsGrade.Range("C13").Value = 4

# New requirement row: "PR sent to lab partner" worth 3 points
$wsGrade.Range("A20").Value = "PR sent to lab partner"
$wsGrade.Range("B20").Value = 3
$wsGrade.Range("C20").Value = 3

# "Site running on Azure": Possible/Actual 5 -> 4
$wsGrade.Range("B21").Value = 4
$wsGrade.Range("C21").Value = 4

# --- View state --------------------------------------------------------
# Leave the Grade sheet's own selection where the author left it, then
# switch back to (and finish on) the Rubric tab, which becomes active.
$wsGrade.Activate()
$wsGrade.Range("A10:B21").Select()

$wsRubric.Activate()
$wsRubric.Range("D13").Select()
